$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet view: scroll position and selected cell
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("O23").Select()

# Add subtotal formulas (Points / Max Points / Difference) for each feature block
$ws.Range("E21").Formula = "=SUM(B21:B26)"
$ws.Range("F21").Formula = "=SUM(C21:C26)"
$ws.Range("G21").Formula = "=F21-E21"

$ws.Range("E29").Formula = "=SUM(B29:B34)"
$ws.Range("F29").Formula = "=SUM(C29:C34)"
$ws.Range("G29").Formula = "=F29-E29"

$ws.Range("E36").Formula = "=SUM(B36)"
$ws.Range("F36").Formula = "=SUM(C36)"
$ws.Range("G36").Formula = "=F36-E36"
# Row 36 carries an explicit row-level style (customFormat), which Excel
# would otherwise stamp onto any newly-touched cell in that row; match the
# un-styled look of the other new helper cells by re-applying the workbook
# default ("Standard") style picked up from a plain cell in the same row.
$ws.Range("E36:G36").Style = $ws.Range("B36").Style

$ws.Range("E39").Formula = "=SUM(B39:B42)"
$ws.Range("F39").Formula = "=SUM(C39:C42)"
$ws.Range("G39").Formula = "=F39-E39"

$ws.Range("E47").Formula = "=SUM(B47:B52)"
$ws.Range("F47").Formula = "=SUM(C47:C52)"
$ws.Range("G47").Formula = "=F47-E47"
